$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value2 = 18.9570735570164
$ws.Cells.Item(2, 2).Value2 = 14.34218387577806
$ws.Cells.Item(2, 3).Value2 = 23.62340838986937
$ws.Cells.Item(3, 1).Value2 = 24.05103843076496
$ws.Cells.Item(3, 2).Value2 = 15.22006587587772
$ws.Cells.Item(3, 3).Value2 = 32.96868354207432
$ws.Cells.Item(4, 1).Value2 = 23.66741125480683
$ws.Cells.Item(4, 2).Value2 = 19.95959585435053
$ws.Cells.Item(4, 3).Value2 = 27.8654980739039
$ws.Cells.Item(5, 1).Value2 = 29.17863098445781
$ws.Cells.Item(5, 2).Value2 = 20.68732493380909
$ws.Cells.Item(5, 3).Value2 = 38.89319879509092
$ws.Cells.Item(6, 1).Value2 = 29.00513344113311
$ws.Cells.Item(6, 2).Value2 = 25.74425588964452
$ws.Cells.Item(6, 3).Value2 = 32.51223808671732
$ws.Cells.Item(7, 1).Value2 = 17.50493707208367
$ws.Cells.Item(7, 2).Value2 = 14.0570865368999
$ws.Cells.Item(7, 3).Value2 = 21.28645448559294
$ws.Cells.Item(8, 1).Value2 = 5.455136183207044
$ws.Cells.Item(8, 2).Value2 = 4.148587976937096
$ws.Cells.Item(8, 3).Value2 = 7.144987446815967
$ws.Cells.Item(9, 1).Value2 = 26.7049372022754
$ws.Cells.Item(9, 2).Value2 = 23.7683122819085
$ws.Cells.Item(9, 3).Value2 = 29.77421483562581
$ws.Cells.Item(10, 1).Value2 = 5.532831667461014
$ws.Cells.Item(10, 2).Value2 = 4.202621132062809
$ws.Cells.Item(10, 3).Value2 = 7.22514953415899
$ws.Cells.Item(11, 1).Value2 = 13.79796146869829
$ws.Cells.Item(11, 2).Value2 = 10.2248508848721
$ws.Cells.Item(11, 3).Value2 = 17.52881965983106
$ws.Cells.Item(12, 1).Value2 = 26.91952053799127
$ws.Cells.Item(12, 2).Value2 = 21.86186015810158
$ws.Cells.Item(12, 3).Value2 = 32.04697448543714
$ws.Cells.Item(13, 1).Value2 = 7.349631581690655
$ws.Cells.Item(13, 2).Value2 = 5.494577388400234
$ws.Cells.Item(13, 3).Value2 = 9.921235393791981
$ws.Cells.Item(14, 1).Value2 = 9.980989690087547
$ws.Cells.Item(14, 2).Value2 = 7.386329313808751
$ws.Cells.Item(14, 3).Value2 = 13.48695935417248
$ws.Cells.Item(15, 1).Value2 = 34.43260979953207
$ws.Cells.Item(15, 2).Value2 = 24.17416525094167
$ws.Cells.Item(15, 3).Value2 = 44.63777810767615
$ws.Cells.Item(16, 1).Value2 = 28.968904646991
$ws.Cells.Item(16, 2).Value2 = 25.62104916887437
$ws.Cells.Item(16, 3).Value2 = 32.51514477160296
$ws.Cells.Item(17, 1).Value2 = 26.74296139459385
$ws.Cells.Item(17, 2).Value2 = 21.80124644595766
$ws.Cells.Item(17, 3).Value2 = 32.05620547191757
$ws.Cells.Item(18, 1).Value2 = 23.6455447383955
$ws.Cells.Item(18, 2).Value2 = 19.96909236882326
$ws.Cells.Item(18, 3).Value2 = 27.80610676470266
$ws.Cells.Item(19, 1).Value2 = 14.04529609762709
$ws.Cells.Item(19, 2).Value2 = 10.44491932184192
$ws.Cells.Item(19, 3).Value2 = 18.03431651941973
$ws.Cells.Item(20, 1).Value2 = 32.38140706538771
$ws.Cells.Item(20, 2).Value2 = 28.48972388569711
$ws.Cells.Item(20, 3).Value2 = 36.74604413337329
$ws.Cells.Item(21, 1).Value2 = 33.23321737125198
$ws.Cells.Item(21, 2).Value2 = 27.914140291793
$ws.Cells.Item(21, 3).Value2 = 38.39260699930262
$ws.Cells.Item(22, 1).Value2 = 5.53573437452109
$ws.Cells.Item(22, 2).Value2 = 4.226387292941371
$ws.Cells.Item(22, 3).Value2 = 7.318248968211521
$ws.Cells.Item(23, 1).Value2 = 36.99919760106646
$ws.Cells.Item(23, 2).Value2 = 26.07221890718393
$ws.Cells.Item(23, 3).Value2 = 46.90225311422616
$ws.Cells.Item(24, 1).Value2 = 29.016646087778
$ws.Cells.Item(24, 2).Value2 = 25.70322242687318
$ws.Cells.Item(24, 3).Value2 = 32.54009159701216
$ws.Cells.Item(25, 1).Value2 = 24.60923423145547
$ws.Cells.Item(25, 2).Value2 = 22.02487064188486
$ws.Cells.Item(25, 3).Value2 = 27.91356365845965
$ws.Cells.Item(26, 1).Value2 = 32.94949922453805
$ws.Cells.Item(26, 2).Value2 = 28.65498831389147
$ws.Cells.Item(26, 3).Value2 = 37.72416267399473
$ws.Cells.Item(27, 1).Value2 = 7.325916997610744
$ws.Cells.Item(27, 2).Value2 = 5.495153997251776
$ws.Cells.Item(27, 3).Value2 = 9.716220683647405
$ws.Cells.Item(28, 1).Value2 = 18.58424620436234
$ws.Cells.Item(28, 2).Value2 = 14.79535010551144
$ws.Cells.Item(28, 3).Value2 = 22.38552107370545
$ws.Cells.Item(29, 1).Value2 = 38.84990876224569
$ws.Cells.Item(29, 2).Value2 = 27.73870665922277
$ws.Cells.Item(29, 3).Value2 = 48.6304151911552
$ws.Cells.Item(30, 1).Value2 = 25.59880620310578
$ws.Cells.Item(30, 2).Value2 = 21.42629082091741
$ws.Cells.Item(30, 3).Value2 = 30.42463341817272
$ws.Cells.Item(31, 1).Value2 = 9.850402445793954
$ws.Cells.Item(31, 2).Value2 = 7.164645217925788
$ws.Cells.Item(31, 3).Value2 = 12.7908940947955
$ws.Cells.Item(32, 1).Value2 = 7.294181404841294
$ws.Cells.Item(32, 2).Value2 = 5.632348202250451
$ws.Cells.Item(32, 3).Value2 = 9.404015400548429
$ws.Cells.Item(33, 1).Value2 = 9.980260739491495
$ws.Cells.Item(33, 2).Value2 = 7.433048372908535
$ws.Cells.Item(33, 3).Value2 = 13.2862691285902
$ws.Cells.Item(34, 1).Value2 = 20.15665814967981
$ws.Cells.Item(34, 2).Value2 = 15.21095907256237
$ws.Cells.Item(34, 3).Value2 = 25.66260584574433
$ws.Cells.Item(35, 1).Value2 = 30.53963244171263
$ws.Cells.Item(35, 2).Value2 = 26.79128099575513
$ws.Cells.Item(35, 3).Value2 = 34.42056449657188
$ws.Cells.Item(36, 1).Value2 = 17.5159596234457
$ws.Cells.Item(36, 2).Value2 = 14.02468033547414
$ws.Cells.Item(36, 3).Value2 = 21.26269045951177
$ws.Cells.Item(37, 1).Value2 = 35.2387365620642
$ws.Cells.Item(37, 2).Value2 = 30.39178991367649
$ws.Cells.Item(37, 3).Value2 = 40.46562691966284
$ws.Cells.Item(38, 1).Value2 = 9.932704125578306
$ws.Cells.Item(38, 2).Value2 = 7.389125480664432
$ws.Cells.Item(38, 3).Value2 = 13.48260049194242
$ws.Cells.Item(39, 1).Value2 = 27.34339676536694
$ws.Cells.Item(39, 2).Value2 = 22.37472426846752
$ws.Cells.Item(39, 3).Value2 = 33.03423358086433
$ws.Cells.Item(40, 1).Value2 = 28.34186644565036
$ws.Cells.Item(40, 2).Value2 = 24.02237942792653
$ws.Cells.Item(40, 3).Value2 = 32.63086630277689
$ws.Cells.Item(41, 1).Value2 = 16.63493009862139
$ws.Cells.Item(41, 2).Value2 = 13.1253911725614
$ws.Cells.Item(41, 3).Value2 = 20.04702908458824
$ws.Cells.Item(42, 1).Value2 = 14.04562228639018
$ws.Cells.Item(42, 2).Value2 = 10.49992634829437
$ws.Cells.Item(42, 3).Value2 = 18.25342717700766
$ws.Cells.Item(43, 1).Value2 = 25.59588048909535
$ws.Cells.Item(43, 2).Value2 = 21.41470070141937
$ws.Cells.Item(43, 3).Value2 = 30.41058268503878
$ws.Cells.Item(44, 1).Value2 = 12.80781167259965
$ws.Cells.Item(44, 2).Value2 = 9.629953048598122
$ws.Cells.Item(44, 3).Value2 = 16.28562918272094
$ws.Cells.Item(45, 1).Value2 = 14.10647202881113
$ws.Cells.Item(45, 2).Value2 = 10.45176376322086
$ws.Cells.Item(45, 3).Value2 = 19.31737129590746
$ws.Cells.Item(46, 1).Value2 = 9.591780362636088
$ws.Cells.Item(46, 2).Value2 = 7.107565871085508
$ws.Cells.Item(46, 3).Value2 = 12.48195241889388
$ws.Cells.Item(47, 1).Value2 = 5.623224210877397
$ws.Cells.Item(47, 2).Value2 = 4.222332441390872
$ws.Cells.Item(47, 3).Value2 = 7.504735185271901
$ws.Cells.Item(48, 1).Value2 = 26.99852048884541
$ws.Cells.Item(48, 2).Value2 = 21.91032248142294
$ws.Cells.Item(48, 3).Value2 = 32.27854485325339
$ws.Cells.Item(49, 1).Value2 = 26.06190488235943
$ws.Cells.Item(49, 2).Value2 = 22.24453825749503
$ws.Cells.Item(49, 3).Value2 = 30.05347133971882
$ws.Cells.Item(50, 1).Value2 = 35.26322568617744
$ws.Cells.Item(50, 2).Value2 = 30.35525712372663
$ws.Cells.Item(50, 3).Value2 = 40.5252350821463
$ws.Cells.Item(51, 1).Value2 = 24.00431122144721
$ws.Cells.Item(51, 2).Value2 = 15.19903618107514
$ws.Cells.Item(51, 3).Value2 = 32.92252109589037
$ws.Cells.Item(52, 1).Value2 = 9.605612682251341
$ws.Cells.Item(52, 2).Value2 = 7.167633154054244
$ws.Cells.Item(52, 3).Value2 = 12.46478982939704
$ws.Cells.Item(53, 1).Value2 = 27.44492359098563
$ws.Cells.Item(53, 2).Value2 = 22.34333084900578
$ws.Cells.Item(53, 3).Value2 = 33.38529316510497
$ws.Cells.Item(54, 1).Value2 = 34.35848693821551
$ws.Cells.Item(54, 2).Value2 = 24.09308451921043
$ws.Cells.Item(54, 3).Value2 = 44.5549987645003
$ws.Cells.Item(55, 1).Value2 = 33.05184500308701
$ws.Cells.Item(55, 2).Value2 = 28.13005749887061
$ws.Cells.Item(55, 3).Value2 = 37.97572465534491
$ws.Cells.Item(56, 1).Value2 = 9.910116111418729
$ws.Cells.Item(56, 2).Value2 = 7.405741756598503
$ws.Cells.Item(56, 3).Value2 = 12.679689085024
$ws.Cells.Item(57, 1).Value2 = 7.291213763337649
$ws.Cells.Item(57, 2).Value2 = 5.578914139592301
$ws.Cells.Item(57, 3).Value2 = 9.478278011861203
$ws.Cells.Item(58, 1).Value2 = 35.42784347285679
$ws.Cells.Item(58, 2).Value2 = 30.46860549941984
$ws.Cells.Item(58, 3).Value2 = 40.63241468443937
$ws.Cells.Item(59, 1).Value2 = 13.66743520524593
$ws.Cells.Item(59, 2).Value2 = 10.27327891860019
$ws.Cells.Item(59, 3).Value2 = 17.07454446492616
$ws.Cells.Item(60, 1).Value2 = 37.40383679609374
$ws.Cells.Item(60, 2).Value2 = 26.53758588438364
$ws.Cells.Item(60, 3).Value2 = 47.27882000785757
$ws.Cells.Item(61, 1).Value2 = 9.993775162492902
$ws.Cells.Item(61, 2).Value2 = 7.484358607127126
$ws.Cells.Item(61, 3).Value2 = 13.18798833164569
$ws.Cells.Item(62, 1).Value2 = 13.81230336371747
$ws.Cells.Item(62, 2).Value2 = 10.25052019671252
$ws.Cells.Item(62, 3).Value2 = 17.73091068039987
$ws.Cells.Item(63, 1).Value2 = 9.945914846611837
$ws.Cells.Item(63, 2).Value2 = 7.390744243320341
$ws.Cells.Item(63, 3).Value2 = 13.4862561826822
$ws.Cells.Item(64, 1).Value2 = 37.40982056694485
$ws.Cells.Item(64, 2).Value2 = 26.4157506418373
$ws.Cells.Item(64, 3).Value2 = 47.09086408450853
$ws.Cells.Item(65, 1).Value2 = 7.138592359475341
$ws.Cells.Item(65, 2).Value2 = 5.374712361678262
$ws.Cells.Item(65, 3).Value2 = 9.339119789812427
$ws.Cells.Item(66, 1).Value2 = 20.04264981842391
$ws.Cells.Item(66, 2).Value2 = 15.22501579357021
$ws.Cells.Item(66, 3).Value2 = 25.28509606084798
$ws.Cells.Item(67, 1).Value2 = 19.48240768588446
$ws.Cells.Item(67, 2).Value2 = 14.71293002663407
$ws.Cells.Item(67, 3).Value2 = 25.05429680521529
$ws.Cells.Item(68, 1).Value2 = 30.99199890778382
$ws.Cells.Item(68, 2).Value2 = 26.32671197989405
$ws.Cells.Item(68, 3).Value2 = 35.85106936582758
$ws.Cells.Item(69, 1).Value2 = 5.577850634113587
$ws.Cells.Item(69, 2).Value2 = 4.245934687038762
$ws.Cells.Item(69, 3).Value2 = 7.487482489665235
$ws.Cells.Item(70, 1).Value2 = 28.27455178492974
$ws.Cells.Item(70, 2).Value2 = 23.77853751483467
$ws.Cells.Item(70, 3).Value2 = 32.13642582526467
$ws.Cells.Item(71, 1).Value2 = 20.06167327462988
$ws.Cells.Item(71, 2).Value2 = 17.32462388217813
$ws.Cells.Item(71, 3).Value2 = 22.892205612472
$ws.Cells.Item(72, 1).Value2 = 23.15122731679353
$ws.Cells.Item(72, 2).Value2 = 20.06378625055872
$ws.Cells.Item(72, 3).Value2 = 26.85076694407353
